$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the leftover style on A54 (previously an empty styled cell, now real data)
$ws.Cells.Item(54, 1).ClearFormats()

$newRows = @(
    ,@('http://purl.obolibrary.org/obo/OBI_0001404', 'genetic characteristics information', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002087 ', 'containing a specimen function', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002088 ', 'specimen container', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002089 ', 'physical store', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002095 ', 'specimen family creation', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002500 ', 'amniotic fluid specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002501 ', 'bile specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002502 ', 'cerebrospinal fluid specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002503 ', 'feces specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002504 ', 'digestive system fluid or secretion specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002505 ', 'milk specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002506 ', 'pericardial fluid specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002507 ', 'saliva specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002508 ', 'sputum specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002509 ', 'sweat specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002510 ', 'synovial fluid specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002511 ', 'vitreous humor specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002512 ', 'bone marrow specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002513 ', 'placenta specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002514 ', 'peritoneal fluid specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002515 ', 'pleural fluid specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002516 ', 'brain specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002517 ', 'hair specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002518 ', 'prostate gland specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002519 ', 'skeletal muscle tissue specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002520', 'heart specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002521 ', 'renal medulla specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002522 ', 'adrenal gland specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002523 ', 'breast specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002524 ', 'urinary bladder specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002525 ', 'tibial artery specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002526 ', 'skin of body specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002527 ', 'pancreas specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002528 ', 'stomach specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002529 ', 'pituitary gland specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002530 ', 'adipose tissue specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002531 ', 'cortex of kidney specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002532 ', 'esophagus mucosa specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002533 ', 'colon specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002534 ', 'lung specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002535 ', 'esophagus muscularis mucosa specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002536 ', 'cerebral cortex specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002537 ', 'thyroid gland specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002538 ', 'cerebellum specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002539 ', 'tibial nerve specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002540 ', 'coronary artery specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002541 ', 'spleen specimen', 'y')
    ,@('http://purl.obolibrary.org/obo/OBI:0002542 ', 'aorta specimen', 'y')
)

$startRow = 54
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

# Update selection to reflect new used range
$ws.Range("A102").Select()
